$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 17.470401
$ws.Range("H2").Value = 52.411203
$ws.Range("I2").Value = 0.8600988665959021
$ws.Range("J2").Value = 0.8884442399952684
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 109.1447706666667
$ws.Range("N2").Value = 327.434312
$ws.Range("O2").Value = 0.3535542089399963
$ws.Range("P2").Value = 0.3655959674582361
$ws.Range("Q2").Value = 1906.802910599704
$ws.Range("R2").Value = 17161.22619539734
$ws.Range("S2").Value = 0.3040915743895016
$ws.Range("T2").Value = 0.3248116314537675

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 17.470401
$ws.Range("H3").Value = 52.411203
$ws.Range("I3").Value = 0.8600988665959021
$ws.Range("J3").Value = 0.8884442399952684
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 47.980825
$ws.Range("N3").Value = 143.942475
$ws.Range("O3").Value = 0.155424969272891
$ws.Range("P3").Value = 0.1607186127944892
$ws.Range("Q3").Value = 838.244253060825
$ws.Range("R3").Value = 7544.198277547425
$ws.Range("S3").Value = 0.1336808399123164
$ws.Range("T3").Value = 0.1427895257972938

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 17.470401
$ws.Range("H4").Value = 52.411203
$ws.Range("I4").Value = 0.8600988665959021
$ws.Range("J4").Value = 0.8884442399952684
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 62.26741999999999
$ws.Range("N4").Value = 186.80226
$ws.Range("O4").Value = 0.2017037397794264
$ws.Range("P4").Value = 0.2085735992386923
$ws.Range("Q4").Value = 1087.83679663542
$ws.Range("R4").Value = 9790.531169718779
$ws.Range("S4").Value = 0.1734851579724394
$ws.Range("T4").Value = 0.1853060128586977

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 17.470401
$ws.Range("H5").Value = 52.411203
$ws.Range("I5").Value = 0.8600988665959021
$ws.Range("J5").Value = 0.8884442399952684
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 58.81030666666667
$ws.Range("N5").Value = 176.43092
$ws.Range("O5").Value = 0.1905050633580386
$ws.Range("P5").Value = 0.1969935053322898
$ws.Range("Q5").Value = 1027.43964039964
$ws.Range("R5").Value = 9246.956763596761
$ws.Range("S5").Value = 0.1638531890750295
$ws.Range("T5").Value = 0.1750177451289501

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 17.470401
$ws.Range("H6").Value = 52.411203
$ws.Range("I6").Value = 0.8600988665959021
$ws.Range("J6").Value = 0.8884442399952684
$ws.Range("K6").Value = 2
$ws.Range("M6").Value = 30.503993
$ws.Range("N6").Value = 61.007986
$ws.Range("O6").Value = 0.09881201864964768
$ws.Range("P6").Value = 0.06811831517629259
$ws.Range("Q6").Value = 532.9169898111929
$ws.Range("R6").Value = 3197.501938867158
$ws.Range("S6").Value = 0.08498810524661511
$ws.Range("T6").Value = 0.06051932475655943

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.8975426666666667
$ws.Range("H7").Value = 2.692628
$ws.Range("I7").Value = 0.04418761940962108
$ws.Range("J7").Value = 0.04564386429080782
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 109.1447706666667
$ws.Range("N7").Value = 327.434312
$ws.Range("O7").Value = 0.3535542089399963
$ws.Range("P7").Value = 0.3655959674582361
$ws.Range("Q7").Value = 97.96208851688178
$ws.Range("R7").Value = 881.658796651936
$ws.Range("S7").Value = 0.01562271882531021
$ws.Range("T7").Value = 0.01668721272393032

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.8975426666666667
$ws.Range("H8").Value = 2.692628
$ws.Range("I8").Value = 0.04418761940962108
$ws.Range("J8").Value = 0.04564386429080782
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 47.980825
$ws.Range("N8").Value = 143.942475
$ws.Range("O8").Value = 0.155424969272891
$ws.Range("P8").Value = 0.1607186127944892
$ws.Range("Q8").Value = 43.06483761936667
$ws.Range("R8").Value = 387.5835385743
$ws.Range("S8").Value = 0.006867859388982558
$ws.Range("T8").Value = 0.007335818551398554

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.8975426666666667
$ws.Range("H9").Value = 2.692628
$ws.Range("I9").Value = 0.04418761940962108
$ws.Range("J9").Value = 0.04564386429080782
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 62.26741999999999
$ws.Range("N9").Value = 186.80226
$ws.Range("O9").Value = 0.2017037397794264
$ws.Range("P9").Value = 0.2085735992386923
$ws.Range("Q9").Value = 55.88766619325333
$ws.Range("R9").Value = 502.98899573928
$ws.Range("S9").Value = 0.008912808086870543
$ws.Range("T9").Value = 0.00952010505829621

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.8975426666666667
$ws.Range("H10").Value = 2.692628
$ws.Range("I10").Value = 0.04418761940962108
$ws.Range("J10").Value = 0.04564386429080782
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 58.81030666666667
$ws.Range("N10").Value = 176.43092
$ws.Range("O10").Value = 0.1905050633580386
$ws.Range("P10").Value = 0.1969935053322898
$ws.Range("Q10").Value = 52.78475947308445
$ws.Range("R10").Value = 475.0628352577601
$ws.Range("S10").Value = 0.00841796523527076
$ws.Range("T10").Value = 0.008991544823557565

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.8975426666666667
$ws.Range("H11").Value = 2.692628
$ws.Range("I11").Value = 0.04418761940962108
$ws.Range("J11").Value = 0.04564386429080782
$ws.Range("K11").Value = 2
$ws.Range("M11").Value = 30.503993
$ws.Range("N11").Value = 61.007986
$ws.Range("O11").Value = 0.09881201864964768
$ws.Range("P11").Value = 0.06811831517629259
$ws.Range("Q11").Value = 27.37863522120134
$ws.Range("R11").Value = 164.271811327208
$ws.Range("S11").Value = 0.004366267873187012
$ws.Range("T11").Value = 0.003109183133625174

# Row 12
$ws.Range("E12").Value = 2
$ws.Range("G12").Value = 1.944141
$ws.Range("H12").Value = 3.888282
$ws.Range("I12").Value = 0.09571351399447693
$ws.Range("J12").Value = 0.06591189571392365
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 109.1447706666667
$ws.Range("N12").Value = 327.434312
$ws.Range("O12").Value = 0.3535542089399963
$ws.Range("P12").Value = 0.3655959674582361
$ws.Range("Q12").Value = 212.192823588664
$ws.Range("R12").Value = 1273.156941531984
$ws.Range("S12").Value = 0.03383991572518456
$ws.Range("T12").Value = 0.02409712328053828

# Row 13
$ws.Range("E13").Value = 2
$ws.Range("G13").Value = 1.944141
$ws.Range("H13").Value = 3.888282
$ws.Range("I13").Value = 0.09571351399447693
$ws.Range("J13").Value = 0.06591189571392365
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 47.980825
$ws.Range("N13").Value = 143.942475
$ws.Range("O13").Value = 0.155424969272891
$ws.Range("P13").Value = 0.1607186127944892
$ws.Range("Q13").Value = 93.28148909632502
$ws.Range("R13").Value = 559.68893457795
$ws.Range("S13").Value = 0.014876269971592
$ws.Range("T13").Value = 0.01059326844579685

# Row 14
$ws.Range("E14").Value = 2
$ws.Range("G14").Value = 1.944141
$ws.Range("H14").Value = 3.888282
$ws.Range("I14").Value = 0.09571351399447693
$ws.Range("J14").Value = 0.06591189571392365
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 62.26741999999999
$ws.Range("N14").Value = 186.80226
$ws.Range("O14").Value = 0.2017037397794264
$ws.Range("P14").Value = 0.2085735992386923
$ws.Range("Q14").Value = 121.05664418622
$ws.Range("R14").Value = 726.33986511732
$ws.Range("S14").Value = 0.01930577372011646
$ws.Range("T14").Value = 0.0137474813216984

# Row 15
$ws.Range("E15").Value = 2
$ws.Range("G15").Value = 1.944141
$ws.Range("H15").Value = 3.888282
$ws.Range("I15").Value = 0.09571351399447693
$ws.Range("J15").Value = 0.06591189571392365
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 58.81030666666667
$ws.Range("N15").Value = 176.43092
$ws.Range("O15").Value = 0.1905050633580386
$ws.Range("P15").Value = 0.1969935053322898
$ws.Range("Q15").Value = 114.33552841324
$ws.Range("R15").Value = 686.0131704794401
$ws.Range("S15").Value = 0.01823390904773834
$ws.Range("T15").Value = 0.01298421537978215

# Row 16
$ws.Range("E16").Value = 2
$ws.Range("G16").Value = 1.944141
$ws.Range("H16").Value = 3.888282
$ws.Range("I16").Value = 0.09571351399447693
$ws.Range("J16").Value = 0.06591189571392365
$ws.Range("K16").Value = 2
$ws.Range("M16").Value = 30.503993
$ws.Range("N16").Value = 61.007986
$ws.Range("O16").Value = 0.09881201864964768
$ws.Range("P16").Value = 0.06811831517629259
$ws.Range("Q16").Value = 59.30406345501301
$ws.Range("R16").Value = 237.216253820052
$ws.Range("S16").Value = 0.009457645529845569
$ws.Range("T16").Value = 0.00448980728610798
